$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2, A3, A4 with the new consolidated card text
$ws.Range("A2").Value = '(''Nightmare Moon'', [''{4}{B}{B}'', ''Legendary Creature — Alicorn'', ''Flying'', ''As long as it’s nighttime, Nightmare Moon gets +2/+2 and has menace.'', ''{6}: Transform Nightmare Moon. Anypony may activate this ability or help pay the cost. When they do, they become your friend.'', ''6/6'', ''Princess Luna'', ''Legendary Creature — Alicorn'', ''Flying'', ''When this creature transforms into Princess Luna, choose up to six cards you own from outside the game with a moon in their art, then exile those cards. As long as those cards remain exiled, you may cast them, and your friends may cast them with your permission. (Gifts are appreciated.)'', ''4/4''])'
$ws.Range("A3").Value = '(''Princess Twilight Sparkle'', [''{W}{U}'', ''Legendary Creature — Alicorn'', ''Flying'', ''Other Alicorns, Horses, Pegasi, Ponies, and Unicorns you control get +1/+1.'', ''{W}{U}{B}{R}{G}: If you control Applejack, Fluttershy, Pinkie Pie, Rainbow Dash, and Rarity, everypony wins the game.'', ''2/2''])'
$ws.Range("A4").Value = '(''Rarity'', [''{1}{W}{U}'', ''Legendary Creature — Unicorn'', ''Rare and mythic rare spells you cast cost {1} less to cast.'', ''{1}, {T}, Reveal a My Little Pony® toy you own: Until end of turn, another target creature gains protection from each color in that toy’s coat, mane, and outfit.'', ''2/2''])'

# Remove the now-unused rows 5 through 26
$ws.Range("A5:A26").ClearContents() | Out-Null

# Ensure the sheet's used range reflects only A1:A4
$ws.Range("A1").Select() | Out-Null
